# tradexcb_strategy.xlsx - add two sample strategy rows (NIFTY CE / PE legs)
# to the Sheet1 data table, matching the columns defined in row 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original file carried an (empty/no-op) workbook protection stub;
# drop it so the saved workbook no longer advertises protection.
$wb.Unprotect() | Out-Null

# Format the new rows as Text first so values such as "0.0", "21", "1.7"
# and the date "2022-04-07" are stored literally instead of being
# auto-coerced into numbers / dates.
$ws.Range("A2:AE3").NumberFormat = "@"

$row2 = @(
    "Buy","0.0","0.0","MARKET","MIS","5.0","NFO","NIFTY","2022-04-07","NIFTY2240718050CE",
    "0","NO","0","NO","Value","0","Value","0","Value","0",
    "1","NO","new","YES","21","1.7","new","21","NO","new","Default"
)

$row3 = @(
    "Buy","0.0","0.0","MARKET","MIS","5.0","NFO","NIFTY","2022-04-07","NIFTY2240718050PE",
    "0","NO","0","NO","Value","0","Value","0","Value","0",
    "1","NO","new","YES","21","1.7","new","21","NO","new","Default"
)

for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $row2[$i]
    $ws.Cells.Item(3, $i + 1).Value = $row3[$i]
}

# Move the active selection to T2, matching the saved workbook state.
$ws.Range("T2").Select() | Out-Null
